# Estadisticos Segundo Parcial 23 Mayo
# Update statistics (Aprobados/Reprobados/porcentajes/Promedio/Blancos) on the
# "1er Parcial", "2o Parcial" and "Final" sheets to reflect newly graded exams.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("E5").Value = 27
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 90
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 7.4
$ws.Range("E6").Value = 95
$ws.Range("F6").Value = 22
$ws.Range("G6").Value = 81.2
$ws.Range("H6").Value = 18.8
$ws.Range("I6").Value = 7.4
$ws.Range("I9").Value = 7.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("I11").Value = 6.1
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("E16").Value = 257
$ws.Range("F16").Value = 104
$ws.Range("G16").Value = 71.2
$ws.Range("H16").Value = 28.8
$ws.Range("I16").Value = 7.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0

$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("E2").Value = 36
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 7.7
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 7.9
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("E4").Value = 23
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 7.7
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 8.300000000000001
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("E6").Value = 117
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 7.9
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("E7").Value = 28
$ws.Range("F7").Value = 9
$ws.Range("G7").Value = 75.7
$ws.Range("H7").Value = 24.3
$ws.Range("I7").Value = 8.6
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("E8").Value = 38
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 97.40000000000001
$ws.Range("H8").Value = 2.6
$ws.Range("I8").Value = 9.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("E9").Value = 28
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 90.3
$ws.Range("H9").Value = 9.699999999999999
$ws.Range("I9").Value = 8.9
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("E10").Value = 47
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 95.90000000000001
$ws.Range("H10").Value = 4.1
$ws.Range("I10").Value = 9.4
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("E11").Value = 32
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 82.09999999999999
$ws.Range("H11").Value = 17.9
$ws.Range("I11").Value = 8.1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("E12").Value = 33
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 86.8
$ws.Range("H12").Value = 13.2
$ws.Range("I12").Value = 9
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("E13").Value = 206
$ws.Range("F13").Value = 27
$ws.Range("G13").Value = 88.40000000000001
$ws.Range("H13").Value = 11.6
$ws.Range("I13").Value = 8.9
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 81.8
$ws.Range("H14").Value = 18.2
$ws.Range("I14").Value = 8.6
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 81.8
$ws.Range("H15").Value = 18.2
$ws.Range("I15").Value = 8.6
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("E16").Value = 332
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 92
$ws.Range("H16").Value = 8
$ws.Range("I16").Value = 8.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0

$ws = $wb.Worksheets.Item("Final")
$ws.Range("E2").Value = 36
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 7.5
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 8.199999999999999
$ws.Range("E4").Value = 23
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 7.7
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 8.199999999999999
$ws.Range("E6").Value = 117
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 7.9
$ws.Range("E7").Value = 28
$ws.Range("F7").Value = 9
$ws.Range("G7").Value = 75.7
$ws.Range("H7").Value = 24.3
$ws.Range("I7").Value = 8.199999999999999
$ws.Range("E8").Value = 38
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 97.40000000000001
$ws.Range("H8").Value = 2.6
$ws.Range("I8").Value = 9.1
$ws.Range("E9").Value = 28
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 90.3
$ws.Range("H9").Value = 9.699999999999999
$ws.Range("I9").Value = 8.4
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("E10").Value = 47
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 95.90000000000001
$ws.Range("H10").Value = 4.1
$ws.Range("I10").Value = 8.4
$ws.Range("E11").Value = 32
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 82.09999999999999
$ws.Range("H11").Value = 17.9
$ws.Range("I11").Value = 7.2
$ws.Range("E12").Value = 33
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 86.8
$ws.Range("H12").Value = 13.2
$ws.Range("I12").Value = 8.699999999999999
$ws.Range("E13").Value = 206
$ws.Range("F13").Value = 27
$ws.Range("G13").Value = 88.40000000000001
$ws.Range("H13").Value = 11.6
$ws.Range("I13").Value = 8.300000000000001
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("I14").Value = 8.4
$ws.Range("I15").Value = 8.4
$ws.Range("E16").Value = 332
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 92
$ws.Range("H16").Value = 8
$ws.Range("I16").Value = 8.199999999999999
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
